# Update "想去人数" (interest count) figures on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 50
$ws1.Range("F6").Value = 142
$ws1.Range("F7").Value = 350
$ws1.Range("F8").Value = 5006
$ws1.Range("F10").Value = 5250
$ws1.Range("F11").Value = 603
$ws1.Range("F12").Value = 1326
$ws1.Range("F13").Value = 98

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 50
$ws4.Range("F6").Value = 142
$ws4.Range("F8").Value = 350
$ws4.Range("F9").Value = 5006
$ws4.Range("F11").Value = 5250
$ws4.Range("F12").Value = 603
$ws4.Range("F13").Value = 1326
$ws4.Range("F14").Value = 98
